$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Annotation score matrix for rows 2-15, columns E-J
$data = @(
    @(2,2,1,1,2,2),
    @(2,1,1,1,2,2),
    @(2,1,1,1,2,2),
    @(2,1,2,2,2,2),
    @(1,1,1,1,1,1),
    @(1,1,1,1,2,2),
    @(1,1,1,1,1,1),
    @(2,1,1,1,2,2),
    @(2,1,1,1,1,2),
    @(2,1,1,1,2,2),
    @(2,2,1,1,1,2),
    @(2,1,1,2,2,2),
    @(2,1,1,1,2,2),
    @(2,2,1,1,2,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 5
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Sheet view changes: freeze top row, scroll, zoom, selections
$window = $excel.ActiveWindow
$window.SplitColumn = 0
$window.SplitRow = 1
$window.FreezePanes = $true
$window.Zoom = 90
$window.ScrollColumn = 2
$window.ScrollRow = 1

$ws.Range("B1").Select()
$ws.Range("G14").Select()
